$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Build the new "File" query text (Bento object repository revisited):
# this query no longer returns File Type / Breed columns.
$newFileQuery = "MATCH (f:file)-->(parent)`n" +
"WITH DISTINCT f, parent`n" +
"MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
" MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
" MATCH (samp:sample)-->(c) `n" +
" WHERE samp.specific_sample_pathology IN [`"Not Applicable`"]  `n" +
"WITH DISTINCT f, parent, c, demo, diag, s`n" +
"RETURN coalesce(f.file_name, '') AS ``File Name``, `n" +
"        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
"        coalesce(f.file_description, '') AS ``Description``,`n" +
"        coalesce(f.file_format, '') AS ``Format``,`n" +
"        coalesce(f.file_size, '') AS ``Size``,`n" +
"        coalesce(c.case_id, '') AS ``Case ID``, `n" +
"        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
"        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# Replace the B4 cell (FilesTab query) with the revised query text.
$ws.Range("B4").Value = $newFileQuery

# Update the view state: scroll so row 4 is at top and select B4,
# matching the author's last interaction with the sheet before saving.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
